# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview sheet: status text "Ready for handoff" -> "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: fill in "Latest Target File" (w/ hyperlink) and
#    "Latest Handback File" columns, and stamp "Latest Handback DateTime"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"

$ov.Columns.Item(5).ColumnWidth = 29.15
$ov.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("J2").Value = "9a560345-e939-4be9-9cd6-ac9e1cab9775.d812918eb6976de0edf8eb75d2ef0ea838e2bd16.zh-cn.xlf"
$zh.Range("K2").Value = "2016-11-29 06:02:52"
$zh.Range("J3").Value = "a4d86911-3495-41d4-9572-f9128aa1f3bf.0ef897c2d620b9c31f962b213c71a017cbd035de.zh-cn.xlf"
$zh.Range("K3").Value = "2016-11-29 06:02:52"

$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6b2d4d2a95e3b7642003ff2add5863164996289/e2e/9a560345-e939-4be9-9cd6-ac9e1cab9775.md", "", "", "9a560345-e939-4be9-9cd6-ac9e1cab9775.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6b2d4d2a95e3b7642003ff2add5863164996289/e2e/a4d86911-3495-41d4-9572-f9128aa1f3bf.md", "", "", "a4d86911-3495-41d4-9572-f9128aa1f3bf.md")

$zh.Columns.Item(3).ColumnWidth = 29.15
$zh.Columns.Item(9).ColumnWidth = 39.15
$zh.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("J2").Value = "9a560345-e939-4be9-9cd6-ac9e1cab9775.d812918eb6976de0edf8eb75d2ef0ea838e2bd16.de-de.xlf"
$de.Range("K2").Value = "2016-11-29 06:03:11"
$de.Range("J3").Value = "a4d86911-3495-41d4-9572-f9128aa1f3bf.0ef897c2d620b9c31f962b213c71a017cbd035de.de-de.xlf"
$de.Range("K3").Value = "2016-11-29 06:03:11"

$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6b2d4d2a95e3b7642003ff2add5863164996289/e2e/9a560345-e939-4be9-9cd6-ac9e1cab9775.md", "", "", "9a560345-e939-4be9-9cd6-ac9e1cab9775.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6b2d4d2a95e3b7642003ff2add5863164996289/e2e/a4d86911-3495-41d4-9572-f9128aa1f3bf.md", "", "", "a4d86911-3495-41d4-9572-f9128aa1f3bf.md")

$de.Columns.Item(3).ColumnWidth = 29.15
$de.Columns.Item(9).ColumnWidth = 39.15
$de.Columns.Item(10).ColumnWidth = 39.15
